$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated "K" (strikeouts) column values (column G), replacing the old
# "Strike#" derived counts with the new per-appearance K totals.
$kValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 2
    6  = 1
    7  = 2
    8  = 3
    9  = 2
    10 = 3
    11 = 0
    12 = 2
    13 = 0
    14 = 1
    15 = 1
    16 = 1
    17 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
